$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# row 43
$ws.Range("H43").Value = 3953.1667
$ws.Range("I43").Value = 2416
$ws.Range("J43").Value = 4465.5557
$ws.Range("K43").Value = 2416
$ws.Range("L43").Value = 4465.5557
$ws.Range("M43").Value = -2347
$ws.Range("N43").Value = -4603.5557
# row 82
$ws.Range("H82").Value = 10007.333
$ws.Range("I82").Value = 10007.333
$ws.Range("K82").Value = 30021.999
$ws.Range("M82").Value = -29615.999
# row 85
$ws.Range("H85").Value = 10007.333
$ws.Range("I85").Value = 10007.333
$ws.Range("K85").Value = 30021.999
$ws.Range("M85").Value = -28617.999
# row 100
$ws.Range("H100").Value = 10039.667
$ws.Range("I100").Value = 10039.667
$ws.Range("K100").Value = 10039.667
$ws.Range("M100").Value = -9498.666999999999
# row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
# row 127
$ws.Range("H127").Value = 1148.5555
$ws.Range("I127").Value = 1209.125
$ws.Range("J127").Value = 664
$ws.Range("K127").Value = 3627.375
$ws.Range("L127").Value = 1992
$ws.Range("M127").Value = 1332.625
$ws.Range("N127").Value = -11912
# row 132
$ws.Range("H132").Value = 2169.0625
$ws.Range("J132").Value = 2660.5454
$ws.Range("L132").Value = 7981.6362
$ws.Range("N132").Value = -13041.6362
# row 137
$ws.Range("H137").Value = 1858.4166
$ws.Range("I137").Value = 1584.4546
$ws.Range("J137").Value = 2288.9285
$ws.Range("K137").Value = 4753.3638
$ws.Range("L137").Value = 6866.7855
$ws.Range("M137").Value = -2203.3638
$ws.Range("N137").Value = -11966.7855

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 5854.885
$ws.Range("I32").Value = 5854.885
$ws.Range("K32").Value = 5854.885
$ws.Range("M32").Value = -5567.885
# row 62
$ws.Range("H62").Value = 89439.2
$ws.Range("J62").Value = 89439.2
$ws.Range("L62").Value = 89439.2
$ws.Range("N62").Value = -90687.2
# row 65
$ws.Range("H65").Value = 89439.2
$ws.Range("J65").Value = 89439.2
$ws.Range("L65").Value = 268317.6
$ws.Range("N65").Value = -274557.6
# row 74
$ws.Range("H74").Value = 2126.5715
$ws.Range("I74").Value = 1664
$ws.Range("J74").Value = 2959.2
$ws.Range("K74").Value = 1664
$ws.Range("L74").Value = 2959.2
$ws.Range("M74").Value = -790
$ws.Range("N74").Value = -4707.2
# row 77
$ws.Range("H77").Value = 2126.5715
$ws.Range("I77").Value = 1664
$ws.Range("J77").Value = 2959.2
$ws.Range("K77").Value = 8320
$ws.Range("L77").Value = 14796
$ws.Range("M77").Value = -3952
$ws.Range("N77").Value = -23532
# row 132
$ws.Range("H132").Value = 2164.7646
$ws.Range("I132").Value = 1761.7693
$ws.Range("J132").Value = 3474.5
$ws.Range("K132").Value = 5285.3079
$ws.Range("L132").Value = 10423.5
$ws.Range("M132").Value = -2755.3079
$ws.Range("N132").Value = -15483.5

$ws = $wb.Worksheets.Item("BSM")
# row 22
$ws.Range("H22").Value = 357.5
$ws.Range("I22").Value = 340
$ws.Range("J22").Value = 375
$ws.Range("K22").Value = 340
$ws.Range("L22").Value = 375
$ws.Range("M22").Value = -167
$ws.Range("N22").Value = -721
# row 40
$ws.Range("H40").Value = 44999
$ws.Range("J40").Value = 44999
$ws.Range("L40").Value = 44999
$ws.Range("N40").Value = -45529
# row 105
$ws.Range("H105").Value = 1200.8125
$ws.Range("I105").Value = 1296.3846
$ws.Range("J105").Value = 786.6667
$ws.Range("K105").Value = 1296.3846
$ws.Range("L105").Value = 786.6667
$ws.Range("M105").Value = 450.6153999999999
$ws.Range("N105").Value = -4280.6667

$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 551.8570999999999
$ws.Range("I16").Value = 494.66666
$ws.Range("K16").Value = 494.66666
$ws.Range("M16").Value = -207.66666
# row 18
$ws.Range("H18").Value = 21799
$ws.Range("J18").Value = 21799
$ws.Range("L18").Value = 21799
$ws.Range("N18").Value = -22259
# row 31
$ws.Range("H31").Value = 3505.5
$ws.Range("I31").Value = 3505.5
$ws.Range("K31").Value = 3505.5
$ws.Range("M31").Value = -3210.5
# row 34
$ws.Range("H34").Value = 3505.5
$ws.Range("I34").Value = 3505.5
$ws.Range("K34").Value = 3505.5
$ws.Range("M34").Value = -3303.5
# row 92
$ws.Range("H92").Value = 42000
$ws.Range("J92").Value = 42000
$ws.Range("L92").Value = 42000
$ws.Range("N92").Value = -46992
# row 111
$ws.Range("H111").Value = 54999
$ws.Range("J111").Value = 54999
$ws.Range("L111").Value = 54999
$ws.Range("N111").Value = -63179
# row 113
$ws.Range("H113").Value = 551.8570999999999
$ws.Range("I113").Value = 494.66666
$ws.Range("K113").Value = 494.66666
$ws.Range("M113").Value = 1675.33334
# row 132
$ws.Range("H132").Value = 1246.1578
$ws.Range("I132").Value = 899.8461
$ws.Range("J132").Value = 1996.5
$ws.Range("K132").Value = 2699.5383
$ws.Range("L132").Value = 5989.5
$ws.Range("M132").Value = -169.5383000000002
$ws.Range("N132").Value = -11049.5

$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 48.142857
$ws.Range("I2").Value = 26.75
$ws.Range("J2").Value = 76.666664
$ws.Range("K2").Value = 160.5
$ws.Range("L2").Value = 459.999984
$ws.Range("M2").Value = -47.5
$ws.Range("N2").Value = -685.999984
# row 81
$ws.Range("H81").Value = 12999.75
$ws.Range("I81").Value = 9499.5
$ws.Range("J81").Value = 16500
$ws.Range("K81").Value = 28498.5
$ws.Range("L81").Value = 49500
$ws.Range("M81").Value = -27375.5
$ws.Range("N81").Value = -51746
# row 84
$ws.Range("H84").Value = 12999.75
$ws.Range("I84").Value = 9499.5
$ws.Range("J84").Value = 16500
$ws.Range("K84").Value = 85495.5
$ws.Range("L84").Value = 148500
$ws.Range("M84").Value = -79879.5
$ws.Range("N84").Value = -159732
# row 93
$ws.Range("H93").Value = 6636.4375
$ws.Range("J93").Value = 6212.2
$ws.Range("L93").Value = 18636.6
$ws.Range("N93").Value = -22380.6
# row 125
$ws.Range("H125").Value = 5033
$ws.Range("J125").Value = 5033
$ws.Range("L125").Value = 15099
$ws.Range("N125").Value = -24939
# row 139
$ws.Range("H139").Value = 1480
$ws.Range("I139").Value = 1480
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 4440
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 700
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# row 43
$ws.Range("H43").Value = 8701.333000000001
$ws.Range("I43").Value = 2441.6
$ws.Range("K43").Value = 2441.6
$ws.Range("M43").Value = -2290.6
# row 80
$ws.Range("H80").Value = 1224.75
$ws.Range("J80").Value = 1250
$ws.Range("L80").Value = 1250
$ws.Range("N80").Value = -3246
# row 83
$ws.Range("H83").Value = 1224.75
$ws.Range("J83").Value = 1250
$ws.Range("L83").Value = 6250
$ws.Range("N83").Value = -16234

$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 5674
$ws.Range("I40").Value = 5898.6665
$ws.Range("K40").Value = 5898.6665
$ws.Range("M40").Value = -5762.6665
# row 46
$ws.Range("H46").Value = 3271.6128
$ws.Range("I46").Value = 739
$ws.Range("K46").Value = 739
$ws.Range("M46").Value = -551
# row 61
$ws.Range("H61").Value = 1197.75
$ws.Range("I61").Value = 1197.75
$ws.Range("K61").Value = 1197.75
$ws.Range("M61").Value = -995.75
# row 113
$ws.Range("H113").Value = 1197.75
$ws.Range("I113").Value = 1197.75
$ws.Range("K113").Value = 1197.75
$ws.Range("M113").Value = 972.25
# row 136
$ws.Range("H136").Value = 90924880
$ws.Range("I136").Value = 17855
$ws.Range("K136").Value = 53565
$ws.Range("M136").Value = -51015

$ws = $wb.Worksheets.Item("WVR")
# row 70
$ws.Range("H70").Value = 12001.25
$ws.Range("J70").Value = 12001.25
$ws.Range("L70").Value = 12001.25
$ws.Range("N70").Value = -12631.25
# row 73
$ws.Range("H73").Value = 12001.25
$ws.Range("J73").Value = 12001.25
$ws.Range("L73").Value = 12001.25
$ws.Range("N73").Value = -14185.25
# row 107
$ws.Range("H107").Value = 2449
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 2449
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 7347
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -11187
# row 122
$ws.Range("H122").Value = 3947.8333
$ws.Range("I122").Value = 4537.4
$ws.Range("K122").Value = 13612.2
$ws.Range("M122").Value = -11162.2
# row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
